# spring 23 week 10 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$data = @(
    @(2, 14, 3, 6),
    @(4, 15, 5, 5),
    @(4, 3, 5, 17),
    @(4, 13, 5, 7),
    @(2, 4, 4, 16),
    @(9, 3, 5, 17),
    @(3, 18, 4, 2),
    @(5, 16, 8, 4),
    @(3, 3, 2, 17),
    @(4, 12, 6, 8),
    @(4, 6, 2, 14),
    @(3, 7, 2, 13),
    @(8, 16, 7, 4),
    @(1, 16, 3, 4)
)

$startRow = 828
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
    $ws.Range("D$row").Value = $vals[3]
}

# Move the view/selection down to the newly added rows, matching the
# author's on-screen state after entering this week's data.
$excel.ActiveWindow.ScrollRow = 830
$ws.Range("A842").Select()
